$d = $word.ActiveDocument

$pairs = @(
    @("406×4=", "567×2="),
    @("354×4=", "284×7="),
    @("773×7=", "281×2="),
    @("458×5=", "800×6="),
    @("358×4=", "782×4="),
    @("269×9=", "724×2="),
    @("244×6=", "603×5="),
    @("852×8=", "569×4="),
    @("250×6=", "884×8="),
    @("606×7=", "747×5="),
    @("602×9=", "807×2="),
    @("232×2=", "167×3="),
    @("779×8=", "700×8="),
    @("920×8=", "730×9="),
    @("932×8=", "855×6="),
    @("418×6=", "661×3="),
    @("952×6=", "835×8="),
    @("646×7=", "109×7="),
    @("124×5=", "885×2="),
    @("596×7=", "871×2="),
    @("173×5=", "758×8="),
    @("251×7=", "310×8="),
    @("997×4=", "419×9="),
    @("647×3=", "948×6="),
    @("797×5=", "268×7=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
